$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompanyUpdates")

# Duplicate row 2 ("Source - Engagement" / Yes / No) into a new row 3
$ws.Range("A2:C2").Copy()
$ws.Range("A3:C3").PasteSpecial()

# Make CompanyUpdates the active sheet, with the new row 3 selected
$ws.Activate() | Out-Null
$ws.Rows.Item(3).Select() | Out-Null
